$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value to a cell while preserving it as literal TEXT, even
# when the text looks like a number (e.g. "13.22" or "000003"). Plain
# `.Value = "13.22"` gets auto-coerced to a number by Excel, which would
# lose the original (text) cell type / leading zeros. Prefixing with an
# apostrophe forces Excel to store it as text (quotePrefix); ClearFormats
# afterwards drops the number-format/style residue left behind so the cell
# ends up with the default (no explicit style) - matching plain data cells
# elsewhere in this workbook.
# ---------------------------------------------------------------------------
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert the new 2022-Q4 row at the top of the
#    data and shift the existing quarters down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summaryData = @(
    @("2022-Q4", 11, 1.17),
    @("2022-Q3", 9, 1.25),
    @("2022-Q2", 4, 0.6899999999999999),
    @("2022-Q1", 1, 0.06),
    @("2021-Q4", 1, 0.08),
    @("2021-Q3", 4, 0.43),
    @("2021-Q2", 2, 0.33),
    @("2021-Q1", 3, 0.47),
    @("2020-Q4", 3, 0.36)
)

# Row 10 is brand new - copy the formatting of the existing last row (A9)
# onto it first so column A keeps the same style as the rest of the table.
$summary.Cells.Item(9, 1).Copy()
$summary.Cells.Item(10, 1).PasteSpecial(-4122)

for ($i = 0; $i -lt $summaryData.Length; $i++) {
    $row = $i + 2
    $summary.Cells.Item($row, 1).Value = $i
    $summary.Cells.Item($row, 2).Value = $summaryData[$i][0]
    $summary.Cells.Item($row, 3).Value = $summaryData[$i][1]
    $summary.Cells.Item($row, 4).Value = $summaryData[$i][2]
}

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q4" detail sheet right after "总计" (pushing all
#    the other quarter sheets down one position, matching the diff).
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q4"

# Borrow formatting from an existing detail sheet so the new sheet matches
# the look (bold/bordered header row, styled column A) of its siblings.
$template = $wb.Worksheets.Item("2022-Q3")
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2:A12").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    Set-TextValue $newSheet.Cells.Item(1, $i + 2) $headers[$i]
}

# A=index, B=code, C=name, D=scale, E=stockPosition, F=positionShare, G=marketValue, H=rank
$fundRows = @(
    @(159993, "鹏华国证证券龙头ETF", "13.22", "98.15", "3.91", "0.5169", 9),
    @(006682, "景顺长城中证500指数增强A", "17.02", "93.72", "1.89", "0.3217", 5),
    @(000978, "景顺长城量化精选股票", "7.44", "94.11", "2.05", "0.1525", 5),
    @(515760, "华夏中证浙江国资创新发展ETF", "2.14", "99.05", "5.44", "0.1164", 4),
    @(012080, "易方达中证500指数量化增强A", "3.80", "93.92", "0.93", "0.0353", 9),
    @(012081, "易方达中证500指数量化增强C", "1.30", "93.92", "0.93", "0.0121", 9),
    @(015860, "宝盈国证证券龙头指数C", "0.18", "92.92", "3.71", "0.0067", 9),
    @(015859, "宝盈国证证券龙头指数A", "0.14", "92.92", "3.71", "0.0052", 9),
    @(519034, "海富通中证500指数增强A", "0.24", "92.50", "1.68", "0.0040", 2),
    @(009004, "海富通中证500指数增强C", "0.04", "92.50", "1.68", "0.0007", 2),
    @(016935, "景顺长城中证500指数增强C", "0.00", "93.72", "1.89", $null, 5)
)

$codes = @("159993", "006682", "000978", "515760", "012080", "012081", "015860", "015859", "519034", "009004", "016935")

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $row = $i + 2
    $r = $fundRows[$i]
    $newSheet.Cells.Item($row, 1).Value = $i
    Set-TextValue $newSheet.Cells.Item($row, 2) $codes[$i]
    Set-TextValue $newSheet.Cells.Item($row, 3) $r[1]
    Set-TextValue $newSheet.Cells.Item($row, 4) $r[2]
    Set-TextValue $newSheet.Cells.Item($row, 5) $r[3]
    Set-TextValue $newSheet.Cells.Item($row, 6) $r[4]
    if ($i -eq 10) {
        $newSheet.Cells.Item($row, 7).Value = 0
    } else {
        Set-TextValue $newSheet.Cells.Item($row, 7) $r[5]
    }
    $newSheet.Cells.Item($row, 8).Value = $r[6]
}

# ---------------------------------------------------------------------------
# 3. Keep "2020-Q4" (now the last tab) the selected/active sheet, matching
#    the original workbook where the last sheet was the active one.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
